$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 47

$ws.Cells.Item($row, 1).Value = "7Y9H5R"
$ws.Cells.Item($row, 2).Value = "Kit de engranaje de alimentación de papel Epson"
$ws.Cells.Item($row, 3).Value = "TM U220"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 5
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Formula = "=(E47-D47)*G47"
$ws.Cells.Item($row, 9).Formula = "=D47*F47"
$ws.Cells.Item($row, 10).Value = 0
